# Apply the "Updated tests and books" change:
#   1. Swap the tab order of "Criterion 3, Air Speed 0.1" and "Criterion 1, Air Speed 0.1"
#      (this brings along each sheet's own table/column names and data automatically).
#   2. Update the "readme" index sheet: reorder its columns to
#      index, JobNo, Date, Author, sheet_name; bump the Date value; and swap the
#      "Criterion 1"/"Criterion 3" sheet_name entries so they line up with the new
#      tab order.

$wb = $excel.ActiveWorkbook

# --- 1. Swap tab order of the two criterion sheets -------------------------
$critOld4 = $wb.Worksheets.Item("Criterion 3, Air Speed 0.1")
$critOld5 = $wb.Worksheets.Item("Criterion 1, Air Speed 0.1")
$critOld5.Move($critOld4)

# --- 2. Update the readme sheet --------------------------------------------
$ws = $wb.Worksheets.Item("readme")

# Header row: column order becomes index, JobNo, Date, Author, sheet_name
$ws.Range("B1").Value = "JobNo"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Author"
$ws.Range("E1").Value = "sheet_name"

# Row 2 - "Criteria Failing, Air Speed 0.1"
$ws.Range("B2").Value = "/c/e"
$ws.Range("C2").Value = "'20220228"
$ws.Range("D2").Value = "jovyan"
$ws.Range("E2").Value = "Criteria Failing, Air Speed 0.1"

# Row 3 - "Criterion 2, Air Speed 0.1"
$ws.Range("B3").Value = "/c/e"
$ws.Range("C3").Value = "'20220228"
$ws.Range("D3").Value = "jovyan"
$ws.Range("E3").Value = "Criterion 2, Air Speed 0.1"

# Row 4 - now "Criterion 1, Air Speed 0.1" (was "Criterion 3" before the swap)
$ws.Range("B4").Value = "/c/e"
$ws.Range("C4").Value = "'20220228"
$ws.Range("D4").Value = "jovyan"
$ws.Range("E4").Value = "Criterion 1, Air Speed 0.1"

# Row 5 - now "Criterion 3, Air Speed 0.1" (was "Criterion 1" before the swap)
$ws.Range("B5").Value = "/c/e"
$ws.Range("C5").Value = "'20220228"
$ws.Range("D5").Value = "jovyan"
$ws.Range("E5").Value = "Criterion 3, Air Speed 0.1"
